$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.022.00"
$ws.Range("E2").Value = "  +5.56%  "
$ws.Range("D3").Value = "2.600.97"
$ws.Range("E3").Value = "  +6.06%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'588.22"
$ws.Range("E5").Value = "  +3.26%  "
$ws.Range("D6").Value = "'155.31"
$ws.Range("E6").Value = "  +6.19%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").Value = "'0.543"
$ws.Range("E8").Value = "  +2.36%  "
$ws.Range("D9").Value = "2.603.93"
$ws.Range("E9").Value = "  +6.19%  "
$ws.Range("E10").Value = "  +3.75%  "
$ws.Range("D11").Value = "'0.161"
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("E12").Value = "  +3.42%  "
$ws.Range("D13").Value = "'5.31"
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("D14").Value = "'29.58"
$ws.Range("E14").Value = "  +3.20%  "
$ws.Range("D15").Value = "'0.0000184"
$ws.Range("E15").Value = "  +5.30%  "
$ws.Range("D16").Value = "3.063.15"
$ws.Range("E16").Value = "  +5.61%  "
$ws.Range("D17").Value = "65.312.82"
$ws.Range("E17").Value = "  +4.41%  "
$ws.Range("D18").Value = "2.606.86"
$ws.Range("E18").Value = "  +5.84%  "
$ws.Range("D19").Value = "'8.14"
$ws.Range("E19").Value = "  +2.92%  "
$ws.Range("D20").Value = "'11.19"
$ws.Range("E20").Value = "  +3.61%  "
$ws.Range("D21").Value = "'355.14"
$ws.Range("E21").Value = "  +9.78%  "
$ws.Range("E22").Value = "  +5.10%  "
$ws.Range("D23").Value = "'2.27"
$ws.Range("E23").Value = "  +4.99%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'10.16"
$ws.Range("E25").Value = "  +1.70%  "
$ws.Range("D26").Value = "'66.73"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("D27").Value = "'638.99"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("E28").Value = "  +11.51%  "
$ws.Range("D30").Value = "'1.50"
$ws.Range("E30").Value = "  +6.37%  "
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").Value = "'8.25"
$ws.Range("E32").Value = "  +4.62%  "
$ws.Range("D33").Value = "'1.90"
$ws.Range("E33").Value = "  +4.81%  "
$ws.Range("E34").Value = "  +7.01%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.61"
$ws.Range("E35").Value = "  +7.49%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'0.993"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").Value = "'4.98"
$ws.Range("E37").Value = "  +5.61%  "
$ws.Range("E38").Value = "  +8.07%  "
$ws.Range("E39").Value = "  +9.04%  "
$ws.Range("D40").Value = "'19.32"
$ws.Range("E40").Value = "  +4.39%  "
$ws.Range("D41").Value = "'155.01"
$ws.Range("E41").Value = "  +2.82%  "
$ws.Range("D42").Value = "'0.376"
$ws.Range("E42").Value = "  +2.68%  "
$ws.Range("D43").Value = "'1.84"
$ws.Range("E43").Value = "  +7.06%  "
$ws.Range("D44").Value = "'42.06"
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("D45").Value = "'163.80"
$ws.Range("E45").Value = "  +7.66%  "
$ws.Range("D46").Value = "0.0₆0310"
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").Value = "'16.05"
$ws.Range("E48").Value = "  +4.89%  "
$ws.Range("D49").Value = "'3.76"
$ws.Range("E49").Value = "  +6.14%  "
$ws.Range("D50").Value = "'21.81"
$ws.Range("E50").Value = "  +8.58%  "
$ws.Range("E51").Value = "  +5.73%  "
